$wb = $excel.ActiveWorkbook

# ---------- Sheet: Companies ----------
$ws = $wb.Worksheets.Item("Companies")

$ws.Range("A2").Value = 119
$ws.Range("B2").Value = "Bubba"
$ws.Range("C2").Value = "B"
$ws.Range("D2").Value = "www.bubba.com"
$ws.Range("I2").Value = "bubba.jpg"
$ws.Range("J2").Value = "bubbaBD.jpg"
$ws.Range("K2").Value = "bubbaBN.jpg"
$ws.Range("M2").Value = 50
$ws.Range("R2").Value = 53

# ---------- Sheet: Bios ----------
$ws = $wb.Worksheets.Item("Bios")

$ws.Range("A2").Value = 119
$ws.Range("B2").Value = "Name: Bubba Wrestling Federation`nLocation: Bubba City, USA`nFounded: 2020`nSize: Medium`n`nOverview:`nBubba Wrestling Federation (BWF) is a new professional wrestling company that is making waves in the industry. Founded in 2020, BWF has quickly gained a reputation for its exciting matches, talented roster, and innovative storytelling.`n`nRoster:`nBWF boasts a diverse and talented roster of wrestlers from all corners of the globe. From high-flying cruiserweights to powerhouse heavyweights, BWF has something for every wrestling fan. Some of the top stars of BWF include `"The Phenom`" Alex Black, `"The Queen of Chaos`" Ruby Reign, and `"The Canadian Crusher`" Jake Maverick.`n`nShow Format:`nBWF puts on weekly television shows that feature a mix of singles matches, tag team matches, and special events. The company also hosts monthly pay-per-view events that showcase the best of BWF's roster in high-stakes matches and intense rivalries.`n`nTitles:`nBWF currently has four championship titles that are contested for on a regular basis. These titles include the BWF World Heavyweight Championship, the BWF Women's Championship, the BWF Tag Team Championships, and the BWF Cruiserweight Championship. These titles are highly coveted and fiercely contested by the talented wrestlers of BWF.`n`nStorylines:`nBWF is known for its engaging and compelling storytelling. Whether it's a bitter rivalry between two wrestlers or a dramatic betrayal within a tag team, BWF's storylines keep fans on the edge of their seats. The creative team at BWF works tirelessly to craft entertaining and unpredictable storylines that keep fans coming back for more.`n`nProduction:`nBWF spares no expense when it comes to the production of its shows. The company's state-of-the-art production team ensures that each event is visually stunning and professionally executed. From pyrotechnics to elaborate entrances, BWF's production values are top-notch and on par with some of the biggest wrestling companies in the world.`n`nCommunity Involvement:`nBWF is committed to giving back to the community and regularly partners with local charities and organizations to make a positive impact. Whether it's hosting fundraising events or volunteering at local schools, BWF is dedicated to using its platform to support those in need.`n`nOverall, BWF is a dynamic and exciting wrestling company that is quickly making a name for itself in the industry. With its talented roster, engaging storylines, and high-quality production values, BWF is a must-watch for any wrestling fan."

# ---------- Sheet: Notes ----------
$ws = $wb.Worksheets.Item("Notes")

$ws.Range("A2").Value = "Bubba"
$ws.Range("B2").Value = "A new wrestling company"
$ws.Range("D2").Value = "bubba.jpg"
$ws.Range("E2").Value = "bubbaBD.jpg"
$ws.Range("F2").Value = "bubbaBanner.jpg"
$ws.Range("H2").Value = "The logo for 'Bubba' would likely feature bold and modern font in red and black color scheme, with a silhouette of a wrestler throwing a high-flying move."
